# Update countries & provincias Spain
# Applies the data refresh captured in the commit diff:
#  - swap the display order of 4 country pairs whose totals crossed over
#  - refresh the numeric stats for the affected rows
#  - bump the "Datos actualizados" timestamp string

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp header (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Julio de 2020 a las 22:06"

# --- Country name swaps (rank changed, rows keep their position) ---
# Azerbaiyan (row56) <-> Ghana (row57)
$ws.Range("A56").Value = "Ghana"
$ws.Range("A57").Value = "Azerbaiyan"

# Sudan (row77) <-> Costa Rica (row78)
$ws.Range("A77").Value = "Costa Rica"
$ws.Range("A78").Value = "Sudan"

# Nueva Zelanda (row132) <-> Ruanda (row133)
$ws.Range("A132").Value = "Ruanda"
$ws.Range("A133").Value = "Nueva Zelanda"

# San Marino (row154) <-> Angola (row155)
$ws.Range("A154").Value = "Angola"
$ws.Range("A155").Value = "San Marino"

# --- Refreshed statistics ---
# Estados Unidos (row 4)
$ws.Range("B4").Value = 3877501
$ws.Range("C4").Value = 44230
$ws.Range("D4").Value = 1789139
$ws.Range("E4").Value = 1945185
$ws.Range("G4").Value = 300
$ws.Range("H4").Value = 143177

# India (row 6)
$ws.Range("B6").Value = 1118107
$ws.Range("C6").Value = 40243
$ws.Range("D6").Value = 700399
$ws.Range("E6").Value = 390205

# Sudafrica (row 8)
$ws.Range("B8").Value = 364328
$ws.Range("C8").Value = 13449
$ws.Range("D8").Value = 191059
$ws.Range("E8").Value = 168236
$ws.Range("G8").Value = 85
$ws.Range("H8").Value = 5033

# Alemania (row 20)
$ws.Range("B20").Value = 202845
$ws.Range("C20").Value = 273
$ws.Range("E20").Value = 5882

# Ecuador (row 31)
$ws.Range("B31").Value = 74013
$ws.Range("C31").Value = 631
$ws.Range("D31").Value = 31901
$ws.Range("E31").Value = 36799
$ws.Range("G31").Value = 31
$ws.Range("H31").Value = 5313

# Guatemala (row 48)
$ws.Range("B48").Value = 38667
$ws.Range("C48").Value = 625
$ws.Range("E48").Value = 13817
$ws.Range("G48").Value = 36
$ws.Range("H48").Value = 1485

# Ghana (now row 56)
$ws.Range("B56").Value = 27667
$ws.Range("C56").Value = 607
$ws.Range("D56").Value = 23249
$ws.Range("E56").Value = 4270
$ws.Range("G56").Value = 3
$ws.Range("H56").Value = 148

# Azerbaiyan (now row 57)
$ws.Range("B57").Value = 27521
$ws.Range("C57").Value = 388
$ws.Range("D57").Value = 18967
$ws.Range("E57").Value = 8200
$ws.Range("G57").Value = 5
$ws.Range("H57").Value = 354

# Costa de Marfil (row 69)
$ws.Range("B69").Value = 14119
$ws.Range("C69").Value = 207
$ws.Range("D69").Value = 8366
$ws.Range("E69").Value = 5661
$ws.Range("G69").Value = 1
$ws.Range("H69").Value = 92

# Costa Rica (now row 77)
$ws.Range("B77").Value = 11114
$ws.Range("C77").Value = 563
$ws.Range("D77").Value = 2966
$ws.Range("E77").Value = 8086
$ws.Range("G77").Value = 8
$ws.Range("H77").Value = 62

# Sudan (now row 78)
$ws.Range("B78").Value = 10992
$ws.Range("C78").Value = 310
$ws.Range("D78").Value = 5707
$ws.Range("E78").Value = 4592
$ws.Range("G78").Value = 20
$ws.Range("H78").Value = 693

# Suazilandia (row 128)
$ws.Range("B128").Value = 1793
$ws.Range("C128").Value = 64
$ws.Range("D128").Value = 822
$ws.Range("E128").Value = 950

# Ruanda (now row 132)
$ws.Range("B132").Value = 1582
$ws.Range("C132").Value = 43
$ws.Range("D132").Value = 834
$ws.Range("E132").Value = 743
$ws.Range("H132").Value = 5

# Nueva Zelanda (now row 133)
$ws.Range("B133").Value = 1553
$ws.Range("C133").Value = 3
$ws.Range("D133").Value = 1506
$ws.Range("E133").Value = 25
$ws.Range("H133").Value = 22

# Angola (now row 154)
$ws.Range("B154").Value = 705
$ws.Range("C154").Value = 18
$ws.Range("D154").Value = 221
$ws.Range("E154").Value = 455
$ws.Range("H154").Value = 29

# San Marino (now row 155)
$ws.Range("B155").Value = 699
$ws.Range("D155").Value = 656
$ws.Range("E155").Value = 1
$ws.Range("H155").Value = 42
